$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy styles from column O to column P for rows 4-14 so formatting matches
# (column O's cell formats are what the new 2022 column should carry).
$ws.Range("O4:O14").Copy()
$ws.Range("P4:P14").PasteSpecial(-4122)  # xlPasteFormats

# Now fill in the new year column (2022) values.
$ws.Range("P4").Value = 2022

$ws.Range("P5").Value = 1
$ws.Range("P6").Value = "-"
$ws.Range("P7").Value = "-"
$ws.Range("P8").Value = "-"
$ws.Range("P9").Value = "-"
$ws.Range("P10").Value = "-"
$ws.Range("P11").Value = "-"
$ws.Range("P12").Value = 1
$ws.Range("P13").Value = "-"
$ws.Range("P14").Value = "-"

# Update the selection to match target state.
$ws.Range("O21:O22").Select()
